$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) onto the new
# header cells I1:J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I0 (col I) and IF (col J) columns, rows 2-28
$data = @(
    @(8, 9),
    @(9, 9),
    @(8, 8),
    @(5, 6),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(7, 8),
    @(10, 10),
    @(6, 7),
    @(7, 7),
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(6, 6),
    @(8, 8),
    @(5, 5),
    @(10, 10),
    @(3, 3),
    @(6, 6),
    @(6, 6),
    @(2, 3),
    @(5, 5),
    @(5, 5),
    @(4, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
